$p = $ppt.ActivePresentation

# Slide 1: Title slide subtitle text change (shape id=3 -> positional index 2)
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "11 – Jupyter Notebooks"

# Slide 4: Add "Simple pattern recognition" paragraph after "Relative frequency" (shape id=4 -> positional index 3)
$s4 = $p.Slides.Item(4)
$tr = $s4.Shapes.Item(3).TextFrame.TextRange
$para3 = $tr.Paragraphs(3)
$para3.InsertAfter("`rSimple pattern recognition") | Out-Null
